# Update the loan-product identifier text on both sheets (a hyphen was
# inserted after the leading "200"), then leave the workbook with the
# ProductLoanOutput sheet active (cell B1 selected there), and cell B1
# selected on ProductLoanInput as well.

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$newName = "200-MS-EI-DB-DL-REC-NON-RNI-CTRFD-DL-MD-TR-1-ONTIME"

$wsInput.Range("B1").Value  = $newName
$wsOutput.Range("B1").Value = $newName

# Leave the input sheet's last selection on B1.
[void]$wsInput.Range("B1").Select()

# Make the output sheet the active tab, selected at B1.
[void]$wsOutput.Activate()
[void]$wsOutput.Range("B1").Select()
